$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

# --- Helper: write a value that must stay a TEXT cell even when it looks
# like a number (e.g. "0.00", "9.00", "100.00", "0"). Using a leading
# apostrophe forces Excel to store it as text (quotePrefix) and resetting
# the cell Style back to "Normal" afterwards clears the quote-prefix
# formatting flag again, leaving a plain shared-string cell with no style
# override - matching how the original file stores these values.
function Set-TextValue($cell, [string]$text) {
    $needsQuote = $text -match '^[+\-]?[0-9]*\.?[0-9]+$'
    if ($needsQuote) {
        $cell.Value = "'" + $text
        $cell.Style = "Normal"
    } else {
        $cell.Value = $text
    }
}

# --- Helper: write a value that contains an embedded line break. Typing a
# line break into a cell makes Excel auto-fit the row height; calling
# AutoFit right after restores the row to a non-custom height so no stray
# ht=/customHeight= survives on the row.
function Set-MultilineValue($cell, [string]$text, $row) {
    $cell.Value = $text
    $row.AutoFit()
}

# ---------------------------------------------------------------------
# Existing rows: translate placeholder / English text into Polish and
# point a couple of cells at already-existing (but differently indexed)
# shared strings.
# ---------------------------------------------------------------------

Set-TextValue $ws.Cells.Item(19, 6) "Tak"

Set-MultilineValue $ws.Cells.Item(20, 6) "Czy na pewno chcesz `nzresetowac ustawienia? " $ws.Rows.Item(20)

Set-TextValue $ws.Cells.Item(21, 6) "Nie"

Set-TextValue $ws.Cells.Item(33, 6) "<>%"

Set-TextValue $ws.Cells.Item(34, 6) "<>V"

Set-TextValue $ws.Cells.Item(35, 6) "0.00"

Set-TextValue $ws.Cells.Item(36, 6) "0.00"

Set-TextValue $ws.Cells.Item(37, 6) "0.00"

# ---------------------------------------------------------------------
# New rows 38-43: "pojemnosc"/water-calculation related single-use texts.
# ---------------------------------------------------------------------

$ws.Cells.Item(38, 2).Value = "SingleUseId49"
$ws.Cells.Item(38, 3).Value = "Default"
$ws.Cells.Item(38, 4).Value = "Center"
$ws.Cells.Item(38, 5).Value = "LTR"
Set-MultilineValue $ws.Cells.Item(38, 6) "Wysokosc`n butelki" $ws.Rows.Item(38)

$ws.Cells.Item(39, 2).Value = "SingleUseId50"
$ws.Cells.Item(39, 3).Value = "Default"
$ws.Cells.Item(39, 4).Value = "Left"
$ws.Cells.Item(39, 5).Value = "LTR"
$ws.Cells.Item(39, 6).Value = "<value>"

$ws.Cells.Item(40, 2).Value = "SingleUseId51"
$ws.Cells.Item(40, 3).Value = "Default"
$ws.Cells.Item(40, 4).Value = "Left"
$ws.Cells.Item(40, 5).Value = "LTR"
Set-TextValue $ws.Cells.Item(40, 6) "0"

$ws.Cells.Item(41, 2).Value = "SingleUseId52"
$ws.Cells.Item(41, 3).Value = "Default"
$ws.Cells.Item(41, 4).Value = "Left"
$ws.Cells.Item(41, 5).Value = "LTR"
$ws.Cells.Item(41, 6).Value = "cm"

$ws.Cells.Item(42, 2).Value = "SingleUseId53"
$ws.Cells.Item(42, 3).Value = "Default"
$ws.Cells.Item(42, 4).Value = "Left"
$ws.Cells.Item(42, 5).Value = "LTR"
$ws.Cells.Item(42, 6).Value = "<value> ml"

$ws.Cells.Item(43, 2).Value = "SingleUseId54"
$ws.Cells.Item(43, 3).Value = "Default"
$ws.Cells.Item(43, 4).Value = "Left"
$ws.Cells.Item(43, 5).Value = "LTR"
Set-TextValue $ws.Cells.Item(43, 6) "0.00"
